$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    everywhere it appears (Overview E2/F2/E3/F3, zh-cn C2/C3, de-de C2/C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Column widths
# ---------------------------------------------------------------------------
# Overview columns E and F (index 5 and 6) widen
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# zh-cn / de-de column C (index 3) and column K (index 11) widen
$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(11).ColumnWidth = 39.166666666666664

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(11).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 3. zh-cn sheet: Latest Target File (J), Latest Handback File (K),
#    Latest Handback DateTime (L) for rows 2 and 3.
# ---------------------------------------------------------------------------
$zhcnHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnHandbackDate = "2017-01-03 07:14:55"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md"

$zhcn.Range("J2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("J2"), $aMdUrl, "", "", "a.md")
$zhcn.Range("J2").Font.Underline = 2
$zhcn.Range("J2").Font.Color = 15570276
$zhcn.Range("K2").Value = $zhcnHandbackFile
$zhcn.Range("L2").Value = $zhcnHandbackDate

$zhcn.Range("J3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("J3"), $aMdUrl, "", "", "a.md")
$zhcn.Range("J3").Font.Underline = 2
$zhcn.Range("J3").Font.Color = 15570276
$zhcn.Range("K3").Value = $zhcnHandbackFile
$zhcn.Range("L3").Value = $zhcnHandbackDate

# ---------------------------------------------------------------------------
# 4. de-de sheet: same columns, with de-de specific handback file/datetime.
# ---------------------------------------------------------------------------
$dedeHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeHandbackDate = "2017-01-03 07:15:07"

$dede.Range("J2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("J2"), $aMdUrl, "", "", "a.md")
$dede.Range("J2").Font.Underline = 2
$dede.Range("J2").Font.Color = 15570276
$dede.Range("K2").Value = $dedeHandbackFile
$dede.Range("L2").Value = $dedeHandbackDate

$dede.Range("J3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("J3"), $aMdUrl, "", "", "a.md")
$dede.Range("J3").Font.Underline = 2
$dede.Range("J3").Font.Color = 15570276
$dede.Range("K3").Value = $dedeHandbackFile
$dede.Range("L3").Value = $dedeHandbackDate

Write-Output "edit complete"
